# Apply cell updates per diff. NumberFormat is forced to text ("@")
# before assigning values so that numeric-looking strings (e.g. "151.10",
# "7.00", "0.0000242") are preserved exactly as literal text instead of
# being auto-converted to numbers (which would drop trailing zeros or
# switch to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.198.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.966.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.971.39"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.452.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.975.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.958.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.95"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000108"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.55"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.49"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.15"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.60"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "44.22"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +15.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.118"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.288"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.27"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "382.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.91%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.746.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0349"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000218"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.47%  "
